$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Best-effort: reposition/resize the workbook window to match the saved view
# (purely cosmetic; harmless if unsupported by the host).
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 390
    $win.Top = 390
    $win.Width = 17580
    $win.Height = 10845
} catch {}

# Clear the Outlook password values in C3:C16, keeping their cell formatting
$ws.Range("C3:C16").ClearContents()

# Move the active selection to C16, matching the saved view state
$ws.Range("C16").Select()
